# Update the marksheet figures: correct answers and total/max marks.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# "Marking" row - number of correct answers
$ws.Range("B11").Value = 5

# "Total" row - total marks obtained, and "obtained/max" text
$ws.Range("B12").Value = 90
$ws.Range("E12").Value = "90/140"
